# "Instances up to S8" - update the mean-result numbers for this instance run.
$wb = $excel.ActiveWorkbook

# --- Sheet "general": runtime value changed ---
$general = $wb.Sheets.Item("general")
$general.Range("B4").Value = 39.44099998474121

# --- Sheet "x": a few j-indices (column B) were renumbered ---
$x = $wb.Sheets.Item("x")
$x.Range("B8").Value = 12
$x.Range("B10").Value = 13
$x.Range("B13").Value = 11

# --- Sheet "Q": the last three Q values (column C) were rotated ---
$q = $wb.Sheets.Item("Q")
$q.Range("C13").Value = 139.702
$q.Range("C14").Value = 68.25399999999999
$q.Range("C15").Value = 140.052
